# Apply scheduled-task data update to Sheet1:
# - Correct a tiny floating point rounding artifact on A18 (date/time serial)
# - Append a new row (19) with the latest sensor reading

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the rounding of the existing row 18 timestamp value.
$ws.Cells.Item(18, 1).Value = 45866.79191548611

# New row 19 values (2025-07-28 20:00:25 reading).
$ws.Cells.Item(19, 1).Value = 45866.83362572934
$ws.Cells.Item(19, 1).NumberFormat = $ws.Cells.Item(18, 1).NumberFormat

$ws.Cells.Item(19, 2).Value = 2025
$ws.Cells.Item(19, 3).Value = 31
$ws.Cells.Item(19, 4).Value = 14.5
$ws.Cells.Item(19, 5).Value = 85.73999999999999
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 2.29
$ws.Cells.Item(19, 8).Value = "E"
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = "20:00:25"
